# Bug tracker: add a new entry for "connecting to mysql server" issue.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new bug-report row (row 3): No., Issue, Priority, Opened on,
# Opened By, Asignee, Date Resolved, Status.
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "authentication issue while connecting to mysql database"
$ws.Range("D3").Value = "High"
$ws.Range("E3").Value = "25/08/2021"
$ws.Range("F3").Value = "Robert Aldis"
$ws.Range("G3").Value = "Robert Aldis"
$ws.Range("H3").Value = "25/08/2021"
$ws.Range("I3").Value = "resolved"

# Wrap the long issue text so it is fully visible, matching the rest of the
# "Issue" column, and grow the row to fit the wrapped text.
$ws.Range("C3:C16").WrapText = $true
$ws.Rows(3).RowHeight = 42.75
